$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their text formatting so values
# such as "8.12" or "537.31" are not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "58.445.20"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").Value = "3.155.79"
$ws.Range("E3").Value = "  +3.03%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "537.31"
$ws.Range("E5").Value = "  +3.30%  "
$ws.Range("D6").Value = "139.83"
$ws.Range("E6").Value = "  +3.43%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.515"
$ws.Range("E8").Value = "  +9.60%  "
$ws.Range("D9").Value = "7.32"
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("E10").Value = "  +3.59%  "
$ws.Range("E11").Value = "  +5.21%  "
$ws.Range("E12").Value = "  +2.49%  "
$ws.Range("D13").Value = "3.707.94"
$ws.Range("E13").Value = "  +3.32%  "
$ws.Range("D14").Value = "26.11"
$ws.Range("E14").Value = "  +4.34%  "
$ws.Range("E15").Value = "  +6.14%  "
$ws.Range("D16").Value = "58.523.20"
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("D17").Value = "3.179.45"
$ws.Range("E17").Value = "  +3.64%  "
$ws.Range("E18").Value = "  +6.51%  "
$ws.Range("E19").Value = "  +5.56%  "
$ws.Range("E20").Value = "  +6.43%  "
$ws.Range("D21").Value = "377.51"
$ws.Range("E21").Value = "  +8.59%  "
$ws.Range("D22").Value = "5.79"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").Value = "70.18"
$ws.Range("E24").Value = "  +2.15%  "
$ws.Range("E25").Value = "  +4.04%  "
$ws.Range("E26").Value = "  +2.16%  "
$ws.Range("D27").Value = "0.977"
$ws.Range("E27").Value = "  -2.55%  "
$ws.Range("D28").Value = "8.12"
$ws.Range("E28").Value = "  +14.26%  "
$ws.Range("D29").Value = "0.0₃0872"
$ws.Range("E29").Value = "  +3.83%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "1.89"
$ws.Range("E30").Value = "  +2.89%  "
$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").Value = "6.15"
$ws.Range("E31").Value = "  +7.30%  "
$ws.Range("D32").Value = "21.88"
$ws.Range("E32").Value = "  +5.04%  "
$ws.Range("D33").Value = "5.18"
$ws.Range("E33").Value = "  +8.26%  "
$ws.Range("E34").Value = "  +6.07%  "
$ws.Range("D35").Value = "160.88"
$ws.Range("E35").Value = "  +2.18%  "
$ws.Range("E36").Value = "  +4.98%  "
$ws.Range("E37").Value = "  +12.83%  "
$ws.Range("D38").Value = "25.44"
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "2.636.36"
$ws.Range("E39").Value = "  +9.62%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "1.67"
$ws.Range("E40").Value = "  +6.84%  "
$ws.Range("E41").Value = "  +4.81%  "
$ws.Range("D42").Value = "4.21"
$ws.Range("E42").Value = "  +5.35%  "
$ws.Range("E43").Value = "  +6.08%  "
$ws.Range("E44").Value = "  +2.97%  "
$ws.Range("D45").Value = "0.0281"
$ws.Range("E45").Value = "  +8.29%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").Value = "3.200.72"
$ws.Range("E47").Value = "  +3.14%  "
$ws.Range("D48").Value = "0.102"
$ws.Range("E48").Value = "  +11.95%  "
$ws.Range("E49").Value = "  +4.48%  "
$ws.Range("D50").Value = "0.980"
$ws.Range("E50").Value = "  +6.15%  "
$ws.Range("D51").Value = "20.26"
$ws.Range("E51").Value = "  +5.10%  "
